$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1232.1892
$ws.Cells.Item(17, 10).Value = 1232.1892
$ws.Cells.Item(17, 12).Value = 3696.5676
$ws.Cells.Item(17, 14).Value = -4032.5676
$ws.Cells.Item(33, 8).Value = 159.7
$ws.Cells.Item(33, 9).Value = 124.111115
$ws.Cells.Item(33, 10).Value = 480
$ws.Cells.Item(33, 11).Value = 124.111115
$ws.Cells.Item(33, 12).Value = 480
$ws.Cells.Item(33, 13).Value = 104.888885
$ws.Cells.Item(33, 14).Value = -938
$ws.Cells.Item(111, 8).Value = 1626.8572
$ws.Cells.Item(111, 9).Value = 1641.3334
$ws.Cells.Item(111, 10).Value = 1616
$ws.Cells.Item(111, 11).Value = 4924.0002
$ws.Cells.Item(111, 12).Value = 4848
$ws.Cells.Item(111, 13).Value = -1857.0002
$ws.Cells.Item(111, 14).Value = -10982
$ws.Cells.Item(135, 8).Value = 553.8333
$ws.Cells.Item(135, 9).Value = 307.67856
$ws.Cells.Item(135, 11).Value = 2769.10704
$ws.Cells.Item(135, 13).Value = -234.1070399999999
$ws.Cells.Item(137, 8).Value = 1078.0256
$ws.Cells.Item(137, 9).Value = 689.7406999999999
$ws.Cells.Item(137, 10).Value = 1951.6666
$ws.Cells.Item(137, 11).Value = 2069.2221
$ws.Cells.Item(137, 12).Value = 5854.9998
$ws.Cells.Item(137, 13).Value = 480.7779
$ws.Cells.Item(137, 14).Value = -10954.9998
$ws.Cells.Item(141, 8).Value = 568.7692
$ws.Cells.Item(141, 9).Value = 568.7692
$ws.Cells.Item(141, 11).Value = 1706.3076
$ws.Cells.Item(141, 13).Value = 3473.6924

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 6887.9375
$ws.Cells.Item(2, 9).Value = 619.5
$ws.Cells.Item(2, 11).Value = 619.5
$ws.Cells.Item(2, 13).Value = -506.5
$ws.Cells.Item(32, 8).Value = 3977.2131
$ws.Cells.Item(32, 9).Value = 3752.7
$ws.Cells.Item(32, 11).Value = 3752.7
$ws.Cells.Item(32, 13).Value = -3465.7
$ws.Cells.Item(34, 8).Value = 10026.5
$ws.Cells.Item(34, 9).Value = 10025
$ws.Cells.Item(34, 11).Value = 10025
$ws.Cells.Item(34, 13).Value = -9754
$ws.Cells.Item(116, 8).Value = 6887.9375
$ws.Cells.Item(116, 9).Value = 619.5
$ws.Cells.Item(116, 11).Value = 619.5
$ws.Cells.Item(116, 13).Value = 1674.5
$ws.Cells.Item(132, 8).Value = 2225.9583
$ws.Cells.Item(132, 9).Value = 1993.0714
$ws.Cells.Item(132, 10).Value = 2552
$ws.Cells.Item(132, 11).Value = 5979.2142
$ws.Cells.Item(132, 12).Value = 7656
$ws.Cells.Item(132, 13).Value = -3449.2142
$ws.Cells.Item(132, 14).Value = -12716

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 6887.9375
$ws.Cells.Item(3, 9).Value = 619.5
$ws.Cells.Item(3, 11).Value = 619.5
$ws.Cells.Item(3, 13).Value = -505.5
$ws.Cells.Item(20, 8).Value = 2104.8333
$ws.Cells.Item(20, 9).Value = 1824.9
$ws.Cells.Item(20, 10).Value = 3504.5
$ws.Cells.Item(20, 11).Value = 1824.9
$ws.Cells.Item(20, 12).Value = 3504.5
$ws.Cells.Item(20, 13).Value = -1577.9
$ws.Cells.Item(20, 14).Value = -3998.5
$ws.Cells.Item(108, 8).Value = 23999.5
$ws.Cells.Item(108, 10).Value = 23999.5
$ws.Cells.Item(108, 12).Value = 23999.5
$ws.Cells.Item(108, 14).Value = -31679.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 834277.8
$ws.Cells.Item(19, 9).Value = 1250204.1
$ws.Cells.Item(19, 10).Value = 2425.25
$ws.Cells.Item(19, 11).Value = 1250204.1
$ws.Cells.Item(19, 12).Value = 2425.25
$ws.Cells.Item(19, 13).Value = -1250034.1
$ws.Cells.Item(19, 14).Value = -2765.25
$ws.Cells.Item(24, 8).Value = 834277.8
$ws.Cells.Item(24, 9).Value = 1250204.1
$ws.Cells.Item(24, 10).Value = 2425.25
$ws.Cells.Item(24, 11).Value = 1250204.1
$ws.Cells.Item(24, 12).Value = 2425.25
$ws.Cells.Item(24, 13).Value = -1250034.1
$ws.Cells.Item(24, 14).Value = -2765.25
$ws.Cells.Item(31, 8).Value = 1196.1143
$ws.Cells.Item(31, 9).Value = 830.6667
$ws.Cells.Item(31, 11).Value = 830.6667
$ws.Cells.Item(31, 13).Value = -535.6667
$ws.Cells.Item(34, 8).Value = 1196.1143
$ws.Cells.Item(34, 9).Value = 830.6667
$ws.Cells.Item(34, 11).Value = 830.6667
$ws.Cells.Item(34, 13).Value = -628.6667
$ws.Cells.Item(58, 8).Value = 1012.1818
$ws.Cells.Item(58, 9).Value = 979.4286
$ws.Cells.Item(58, 11).Value = 979.4286
$ws.Cells.Item(58, 13).Value = -776.4286
$ws.Cells.Item(122, 8).Value = 975.4545000000001
$ws.Cells.Item(122, 9).Value = 914.44446
$ws.Cells.Item(122, 10).Value = 1250
$ws.Cells.Item(122, 11).Value = 2743.33338
$ws.Cells.Item(122, 12).Value = 3750
$ws.Cells.Item(122, 13).Value = -293.33338
$ws.Cells.Item(122, 14).Value = -8650
$ws.Cells.Item(132, 8).Value = 6730.9546
$ws.Cells.Item(132, 9).Value = 8588.071
$ws.Cells.Item(132, 10).Value = 3481
$ws.Cells.Item(132, 11).Value = 25764.213
$ws.Cells.Item(132, 12).Value = 10443
$ws.Cells.Item(132, 13).Value = -23234.213
$ws.Cells.Item(132, 14).Value = -15503
$ws.Cells.Item(134, 8).Value = 1959.4
$ws.Cells.Item(134, 9).Value = 2035.875
$ws.Cells.Item(134, 10).Value = 1653.5
$ws.Cells.Item(134, 11).Value = 6107.625
$ws.Cells.Item(134, 12).Value = 4960.5
$ws.Cells.Item(134, 13).Value = -3572.625
$ws.Cells.Item(134, 14).Value = -10030.5
$ws.Cells.Item(135, 8).Value = 31000
$ws.Cells.Item(135, 10).Value = 31000
$ws.Cells.Item(135, 12).Value = 31000
$ws.Cells.Item(135, 14).Value = -41140
$ws.Cells.Item(136, 8).Value = 1012.1818
$ws.Cells.Item(136, 9).Value = 979.4286
$ws.Cells.Item(136, 11).Value = 2938.2858
$ws.Cells.Item(136, 13).Value = -388.2857999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1274.0968
$ws.Cells.Item(5, 9).Value = 1336.9259
$ws.Cells.Item(5, 11).Value = 4010.7777
$ws.Cells.Item(5, 13).Value = -3898.7777
$ws.Cells.Item(6, 8).Value = 164.25
$ws.Cells.Item(6, 9).Value = 213.33333
$ws.Cells.Item(6, 10).Value = 17
$ws.Cells.Item(6, 11).Value = 639.99999
$ws.Cells.Item(6, 12).Value = 51
$ws.Cells.Item(6, 13).Value = -526.99999
$ws.Cells.Item(6, 14).Value = -277
$ws.Cells.Item(108, 8).Value = 1348.8889
$ws.Cells.Item(108, 9).Value = 440
$ws.Cells.Item(108, 10).Value = 1803.3334
$ws.Cells.Item(108, 11).Value = 1320
$ws.Cells.Item(108, 12).Value = 5410.0002
$ws.Cells.Item(108, 13).Value = 1560
$ws.Cells.Item(108, 14).Value = -11170.0002
$ws.Cells.Item(112, 8).Value = 13825
$ws.Cells.Item(112, 9).Value = 4000
$ws.Cells.Item(112, 10).Value = 15790
$ws.Cells.Item(112, 11).Value = 12000
$ws.Cells.Item(112, 12).Value = 47370
$ws.Cells.Item(112, 13).Value = -10892
$ws.Cells.Item(112, 14).Value = -49586
$ws.Cells.Item(131, 8).Value = 1223.9691
$ws.Cells.Item(131, 10).Value = 1239.2211
$ws.Cells.Item(131, 12).Value = 3717.6633
$ws.Cells.Item(131, 14).Value = -13797.6633
$ws.Cells.Item(132, 8).Value = 1773.091
$ws.Cells.Item(132, 10).Value = 4500
$ws.Cells.Item(132, 12).Value = 40500
$ws.Cells.Item(132, 14).Value = -45560
$ws.Cells.Item(135, 8).Value = 1274.0968
$ws.Cells.Item(135, 9).Value = 1336.9259
$ws.Cells.Item(135, 11).Value = 12032.3331
$ws.Cells.Item(135, 13).Value = -9497.3331
$ws.Cells.Item(140, 8).Value = 35977.465
$ws.Cells.Item(140, 9).Value = 57971.777
$ws.Cells.Item(140, 11).Value = 173915.331
$ws.Cells.Item(140, 13).Value = -168735.331

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2817.6487
$ws.Cells.Item(102, 9).Value = 2081.3
$ws.Cells.Item(102, 10).Value = 5973.4287
$ws.Cells.Item(102, 11).Value = 2081.3
$ws.Cells.Item(102, 12).Value = 5973.4287
$ws.Cells.Item(102, 13).Value = -459.3000000000002
$ws.Cells.Item(102, 14).Value = -9217.4287
$ws.Cells.Item(126, 8).Value = 2176.6667
$ws.Cells.Item(126, 9).Value = 1765
$ws.Cells.Item(126, 11).Value = 5295
$ws.Cells.Item(126, 13).Value = -2825
$ws.Cells.Item(132, 8).Value = 2627.9656
$ws.Cells.Item(132, 9).Value = 2254.8667
$ws.Cells.Item(132, 11).Value = 6764.6001
$ws.Cells.Item(132, 13).Value = -4234.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2091
$ws.Cells.Item(7, 9).Value = 2043
$ws.Cells.Item(7, 10).Value = 2355
$ws.Cells.Item(7, 11).Value = 2043
$ws.Cells.Item(7, 12).Value = 2355
$ws.Cells.Item(7, 13).Value = -1931
$ws.Cells.Item(7, 14).Value = -2579
$ws.Cells.Item(40, 8).Value = 3644.7083
$ws.Cells.Item(40, 9).Value = 2723.3125
$ws.Cells.Item(40, 11).Value = 2723.3125
$ws.Cells.Item(40, 13).Value = -2587.3125
$ws.Cells.Item(122, 8).Value = 14714177
$ws.Cells.Item(122, 9).Value = 27790586
$ws.Cells.Item(122, 10).Value = 3216.25
$ws.Cells.Item(122, 11).Value = 83371758
$ws.Cells.Item(122, 12).Value = 9648.75
$ws.Cells.Item(122, 13).Value = -83369308
$ws.Cells.Item(122, 14).Value = -14548.75
$ws.Cells.Item(126, 8).Value = 2091
$ws.Cells.Item(126, 9).Value = 2043
$ws.Cells.Item(126, 10).Value = 2355
$ws.Cells.Item(126, 11).Value = 6129
$ws.Cells.Item(126, 12).Value = 7065
$ws.Cells.Item(126, 13).Value = -3659
$ws.Cells.Item(126, 14).Value = -12005
$ws.Cells.Item(132, 8).Value = 22247.94
$ws.Cells.Item(132, 9).Value = 1395.44
$ws.Cells.Item(132, 10).Value = 43969.293
$ws.Cells.Item(132, 11).Value = 4186.32
$ws.Cells.Item(132, 12).Value = 131907.879
$ws.Cells.Item(132, 13).Value = -1656.32
$ws.Cells.Item(132, 14).Value = -136967.879
$ws.Cells.Item(136, 8).Value = 6490.476
$ws.Cells.Item(136, 9).Value = 7465
$ws.Cells.Item(136, 10).Value = 2348.75
$ws.Cells.Item(136, 11).Value = 22395
$ws.Cells.Item(136, 12).Value = 7046.25
$ws.Cells.Item(136, 13).Value = -19845
$ws.Cells.Item(136, 14).Value = -12146.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 20001130
$ws.Cells.Item(122, 9).Value = 26001160
$ws.Cells.Item(122, 11).Value = 78003480
$ws.Cells.Item(122, 13).Value = -78001030
$ws.Cells.Item(132, 8).Value = 3166
$ws.Cells.Item(132, 9).Value = 3742.9565
$ws.Cells.Item(132, 10).Value = 2336.625
$ws.Cells.Item(132, 11).Value = 11228.8695
$ws.Cells.Item(132, 12).Value = 7009.875
$ws.Cells.Item(132, 13).Value = -8698.869499999999
$ws.Cells.Item(132, 14).Value = -12069.875
$ws.Cells.Item(136, 8).Value = 626.5454999999999
$ws.Cells.Item(136, 9).Value = 435.78946
$ws.Cells.Item(136, 10).Value = 1834.6666
$ws.Cells.Item(136, 11).Value = 1307.36838
$ws.Cells.Item(136, 12).Value = 5503.9998
$ws.Cells.Item(136, 13).Value = 1242.63162
$ws.Cells.Item(136, 14).Value = -10603.9998

